# Ajout d'un groupe dans les donnees d'exemples de la paroisse de Morges
# pour raison de consistance : insere "Groupes libres" (403060100) comme
# nouveau sous-groupe de 403060000, et reparente les anciens enfants directs
# 403060102..403060132 sous ce nouveau noeud.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the old row 97, shifting everything below down by one.
$ws.Rows(97).Insert()

# Populate the new row with the "Groupes libres" entry.
$ws.Range("A97").Value = 403060100
$ws.Range("B97").Value = "Groupes libres"
$ws.Range("C97").Value = 403060000
$ws.Range("D97").Value = 2010000000

# Re-parent the 31 rows that used to be direct children of 403060000
# (old rows 97..127, now rows 98..128) under the new 403060100 group.
$ws.Range("C98:C128").Value = 403060100

# Restore the view state roughly matching the author's saved selection.
$ws.Range("C94").Select()
